$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout: A dispositivo | B sesion_id | C dia | D mes | E año | F timestamp
#                 | G ubicacion | H modelo | I variable | J valor
# Target layout:   A patente | B sesion_id | C dia | D mes | E año | F timestamp
#                 | G variable | H valor

# 1) Drop "dispositivo" (col A)
$ws.Columns.Item(1).Delete()
# columns now: A sesion_id | B dia | C mes | D año | E timestamp | F ubicacion | G modelo | H variable | I valor

# 2) Drop "ubicacion" (col F)
$ws.Columns.Item(6).Delete()
# columns now: A sesion_id | B dia | C mes | D año | E timestamp | F modelo | G variable | H valor

# 3) Drop "modelo" (col F)
$ws.Columns.Item(6).Delete()
# columns now: A sesion_id | B dia | C mes | D año | E timestamp | F variable | G valor

Write-Host "After drops:" $ws.UsedRange.Address()

# 4) Insert a new column at the front for "patente"
$ws.Columns.Item(1).Insert()
Write-Host "After insert:" $ws.UsedRange.Address()

$ws.Range("A1").Value = "patente"

# 5) Replace the "variable" column (G, formerly "% co2") with "°C" readings
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 7).Value = "°C"
}

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = "MP-01-EXPRESS"
}

$null = $ws.Range("K20").Select()
